# Update cryptos list values (price + volume(1h) change); a few rows also
# had the coin name/link change because two entries swapped positions in
# the underlying feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $B, $C, $D, $E) {
    if ($B -ne $null) { $ws.Cells.Item($Row, 2).Value = $B }
    if ($C -ne $null) { $ws.Cells.Item($Row, 3).Value = $C }
    if ($D -ne $null) { $ws.Cells.Item($Row, 4).Value = $D }
    if ($E -ne $null) { $ws.Cells.Item($Row, 5).Value = $E }
}

# Row 2 - Bitcoin
Set-Row 2 $null $null "66.612.74" "  -4.71%  "
# Row 3 - Ethereum
Set-Row 3 $null $null "3.478.48" "  -5.86%  "
# Row 4 - TetherUSD
Set-Row 4 $null $null $null "  -0.22%  "
# Row 5 - BNB
Set-Row 5 $null $null "600.82" "  -7.33%  "
# Row 6 - Solana
Set-Row 6 $null $null "149.07" "  -7.90%  "
# Row 7 - LidoStakedEther
Set-Row 7 $null $null "3.478.30" "  -5.83%  "
# Row 8 - USDC
Set-Row 8 $null $null $null "  -0.08%  "
# Row 9 - XRP
Set-Row 9 $null $null "0.480" "  -4.77%  "
# Row 10 - Dogecoin
Set-Row 10 $null $null "0.138" "  -5.42%  "
# Row 11 - Toncoin
Set-Row 11 $null $null "6.87" "  -4.46%  "
# Row 12 - Cardano
Set-Row 12 $null $null "0.422" "  -5.27%  "
# Row 13 - ShibaInu
Set-Row 13 $null $null "0.0000219" "  -6.15%  "
# Row 14 - WrappedliquidstakedEther2.0
Set-Row 14 $null $null "4.066.68" "  -5.91%  "
# Row 15 - Avalanche
Set-Row 15 $null $null "31.37" "  -4.31%  "
# Row 16 - WrappedEther
Set-Row 16 $null $null "3.478.26" "  -6.03%  "
# Row 17 - WrappedBTC
Set-Row 17 $null $null "66.563.69" "  -4.81%  "
# Row 18 - TRON
Set-Row 18 $null $null $null "  -0.18%  "
# Row 19 - Polkadot
Set-Row 19 $null $null "6.36" "  -2.44%  "
# Row 20 - Chainlink
Set-Row 20 $null $null "15.01" "  -6.66%  "
# Row 21 - BitcoinCash
Set-Row 21 $null $null "442.29" "  -6.34%  "
# Row 22 - Uniswap
Set-Row 22 $null $null "9.03" "  -14.80%  "
# Row 23 - Polygon
Set-Row 23 $null $null "0.623" "  -4.58%  "
# Row 24 - Litecoin
Set-Row 24 $null $null "76.97" "  -4.01%  "
# Row 25 - Dai
Set-Row 25 $null $null "1.00" "  -0.07%  "
# Row 26 - WrappedeETH
Set-Row 26 $null $null "3.608.82" "  -6.14%  "
# Row 27 - PEPE
Set-Row 27 $null $null "0.0000124" "  -3.21%  "
# Row 28 - InternetComputer(DFINITY)
Set-Row 28 $null $null "10.10" "  -8.48%  "
# Row 29 - RenderToken
Set-Row 29 $null $null "8.19" "  -10.65%  "
# Row 30 - PancakeSwap
Set-Row 30 $null $null "2.50" "  -6.16%  "
# Row 31 - Fetch.AI
Set-Row 31 $null $null "1.57" "  -8.95%  "
# Row 32 - Binance-PegBSC-USD
Set-Row 32 $null $null "0.999" "  -0.08%  "

# Row 33 / Row 34 - EthereumClassic and Kaspa swap places
Set-Row 33 "Kaspa" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas" "0.158" "  -4.09%  "
Set-Row 34 "EthereumClassic" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" "25.53" "  -5.19%  "

# Row 35 - NEARProtocol
Set-Row 35 $null $null "6.14" "  -6.59%  "
# Row 36 - ImmutableX
Set-Row 36 $null $null "1.85" "  -8.21%  "
# Row 37 - RenzoRestakedETH
Set-Row 37 $null $null "3.463.67" "  -6.21%  "
# Row 38 - Aptos
Set-Row 38 $null $null "7.94" "  -6.03%  "
# Row 39 - USDe (unchanged)
# Row 40 - FirstDigitalUSD
Set-Row 40 $null $null "0.998" "  -0.33%  "
# Row 41 - Monero
Set-Row 41 $null $null "172.57" "  -4.56%  "
# Row 42 - Stacks
Set-Row 42 $null $null "2.16" "  -4.53%  "

# Row 43 / Row 44 - Filecoin and Hedera swap places
Set-Row 43 "Hedera" "https://coinranking.com/coin/jad286TjB+hedera-hbar" "0.0861" "  -5.26%  "
Set-Row 44 "Filecoin" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" "5.47" "  -7.45%  "

# Row 45 - Mantle
Set-Row 45 $null $null "0.881" "  -5.52%  "
# Row 46 - OKB
Set-Row 46 $null $null "45.11" "  -3.95%  "
# Row 47 - InjectiveProtocol
Set-Row 47 $null $null "26.75" "  -8.35%  "
# Row 48 - ONDO
Set-Row 48 $null $null "1.21" "  -3.82%  "

# Row 49 / Row 50 - Cosmos and dogwifhat swap places
Set-Row 49 "dogwifhat" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif" "2.48" "  -12.21%  "
Set-Row 50 "Cosmos" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom" "7.54" "  -4.07%  "

# Row 51 - SuiNetwork
Set-Row 51 $null $null "1.00" "  -5.44%  "
